$wb = $excel.ActiveWorkbook

# --- Update header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet, placed after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45564.99999999999, 368, 62.36501325868949, 689.657219700605),
    @(45571.99999999999, 327, 32.57262078994179, 632.8603575657144),
    @(45599.99999999999, 162, -139.0981708716037, 450.9992335437506),
    @(45634.99999999999, 0, -346.1353921898136, 245.9528539410119),
    @(45641.99999999999, 0, -375.4709299920406, 227.2299804352228),
    @(45648.99999999999, 0, -414.6694559256749, 189.4348241983309),
    @(45655.99999999999, 0, -462.3765935763157, 128.2034982759944),
    @(45662.99999999999, 0, -516.7439597381905, 82.63089335110094),
    @(45669.99999999999, 0, -536.9554356654284, 90.82313794749193),
    @(45676.99999999999, 0, -605.1695583568827, 18.19779409083752),
    @(45683.99999999999, 0, -636.0852232660031, -35.59555275904577),
    @(45690.99999999999, 0, -684.2399997670816, -67.31326396914487)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# --- Formatting: reuse the existing header / date styles from "Weekly Quantity"
# so the new sheet matches the look of the other two (bold+border header row,
# date-formatted first column), without clobbering the values we just wrote. ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)

$excel.CutCopyMode = $false
